# SRS_Review.xlsx update - "update CYRS review sheet"
#
# Changes:
#  Introduction sheet:
#   - Ref Version (D7): 1 -> 1.1
#   - Last update date (D9): text "27/01/2020" -> real date 02/07/2020
#   - Revision history: new row 14 entry (Version 0.2, T.Sharaby, 02/07/2020, "Add some pints ")
#  Cross review points sheet:
#   - New review rows 9 and 10 filled in with open review points for V1.1

$wb = $excel.ActiveWorkbook
$wsIntro = $wb.Worksheets.Item("Introduction ")
$wsCross = $wb.Worksheets.Item("Cross review points ")

$newDate = Get-Date -Year 2020 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0

# --- Cross review points sheet: add the two new open review points first so
#     the new shared strings land in the same order as the target workbook
#     (long descriptions first, then the "V1.1" detection-version tag, then
#     the Introduction-sheet "Add some pints " note).
$wsCross.Range("F9").Value = "At the first page it says the version is 1.0 and it is proposed , but the in history it is 1.1 "
$wsCross.Range("F10").Value = "Still not all rebiew points resolved "

$wsCross.Range("A9").Value = $newDate
$wsCross.Range("A9").NumberFormatLocal = $wsCross.Range("A2").NumberFormatLocal
$wsCross.Range("B9").Value = "TSH"
$wsCross.Range("C9").Value = "V1.1"
$wsCross.Range("D9").Value = "SRS"
$wsCross.Range("E9").Value = "Software requirement "
$wsCross.Range("H9").Value = "Open"

$wsCross.Range("A10").Value = $newDate
$wsCross.Range("A10").NumberFormatLocal = $wsCross.Range("A2").NumberFormatLocal
$wsCross.Range("B10").Value = "TSH"
$wsCross.Range("C10").Value = "V1.1"
$wsCross.Range("D10").Value = "SRS"
$wsCross.Range("E10").Value = "Software requirement "
$wsCross.Range("H10").Value = "Open"

# --- Introduction sheet: status block ---
# Ref Version 1 -> 1.1
$wsIntro.Range("D7").Value = 1.1

# Last update date: was stored as text "27/01/2020"; replace with a real date
$wsIntro.Range("D9").Value = $newDate
$wsIntro.Range("D9").NumberFormatLocal = "mm-dd-yy"

# --- Introduction sheet: revision history new row (version 0.2) ---
$wsIntro.Range("B14").Value = 0.2
$wsIntro.Range("C14").Value = "T.Sharaby"

# Reuse exactly the same date-number-format style as D9 above (copy the
# format only, so the workbook doesn't end up with two near-identical
# cellXfs entries for the same format).
$wsIntro.Range("D9").Copy()
$wsIntro.Range("E14").PasteSpecial(-4122)
$wsIntro.Range("E14").Value = $newDate

$wsIntro.Range("G14").Value = "Add some pints "

# Restore the selection that Copy() leaves behind
$wsIntro.Range("D8:H8").Select()
